# Spring1415.xlsx: "Added logic for inserting data into Vars table"
#
# The BasicData sheet previously carried a Year/Quarter pair of rows
# (Year=2015, Quarter=Spring) ahead of the Layout_* variables. Those two
# rows are removed (shifting the Layout_* rows up), the sheet is renamed
# from "BasicData" to "Variables", and the active selection moves to N16.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the "Year" (row 2) and "Quarter" (row 3) rows entirely; this
# shifts Layout_Section1.. up to rows 2-6.
$ws.Rows("2:3").Delete()

# Rename the sheet to reflect its new purpose.
$ws.Name = "Variables"

# Move the selection to match the committed state.
$ws.Range("N16").Select()
